# Corrected error in currency conversion factor.
#
# The "About" sheet holds the raw USD<->HKD conversion rate in A35 (with a
# label in B35) that feeds A36 (=A26*A35), which in turn is consumed by the
# OCCF-DpLOCU / OCCF-DpMOCU / OCCF-DpSOCU sheets. The old rate (7.8285,
# labeled "USD converted to HKD") was wrong; replace it with the correct
# rate (~0.1277, i.e. "HKD converted to USD") and fix the label.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Corrected conversion factor and its (now-swapped) label.
$ws.Range("A35").Value = 0.12773839177364757
$ws.Range("B35").Value = "HKD converted to USD"

# The "this number was used" note now sits on row 36 instead of row 35.
$ws.Range("B36").Value = "this number was used"

$wb.Application.Calculate()
